$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 40862.8
$ws.Range("I31").Value = 40862.8
$ws.Range("K31").Value = 122588.4
$ws.Range("M31").Value = -122358.4

$ws.Range("H69").Value = 46804.285
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 54105
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 162315
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -164063

$ws.Range("H72").Value = 46804.285
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 54105
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 486945
$ws.Range("M72").Value = -22632
$ws.Range("N72").Value = -495681

$ws.Range("H106").Value = 4124.5293
$ws.Range("I106").Value = 4182.375
$ws.Range("K106").Value = 4182.375
$ws.Range("M106").Value = -3551.375

$ws.Range("H111").Value = 5955.3125
$ws.Range("I111").Value = 4592
$ws.Range("J111").Value = 8954.6
$ws.Range("K111").Value = 13776
$ws.Range("L111").Value = 26863.8
$ws.Range("M111").Value = -10709
$ws.Range("N111").Value = -32997.8

$ws.Range("H113").Value = 7821
$ws.Range("I113").Value = 5315.6665
$ws.Range("J113").Value = 9700
$ws.Range("K113").Value = 5315.6665
$ws.Range("L113").Value = 9700
$ws.Range("M113").Value = -2061.6665
$ws.Range("N113").Value = -16208

$ws.Range("H132").Value = 14477.308
$ws.Range("I132").Value = 1021.25
$ws.Range("K132").Value = 3063.75
$ws.Range("M132").Value = -533.75

$ws.Range("H141").Value = 3064.8572
$ws.Range("I141").Value = 3064.8572
$ws.Range("K141").Value = 9194.571599999999
$ws.Range("M141").Value = -4014.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3107.5
$ws.Range("I61").Value = 1425.85
$ws.Range("K61").Value = 1425.85
$ws.Range("M61").Value = -1213.85

$ws.Range("H122").Value = 3414.3171
$ws.Range("I122").Value = 2114.12
$ws.Range("K122").Value = 6342.36
$ws.Range("M122").Value = -3892.36

$ws.Range("H136").Value = 3107.5
$ws.Range("I136").Value = 1425.85
$ws.Range("K136").Value = 4277.549999999999
$ws.Range("M136").Value = -1727.549999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 70387.5
$ws.Range("J74").Value = 98000
$ws.Range("L74").Value = 98000
$ws.Range("N74").Value = -99872

$ws.Range("H77").Value = 70387.5
$ws.Range("J77").Value = 98000
$ws.Range("L77").Value = 294000
$ws.Range("N77").Value = -303360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 145.44444
$ws.Range("I7").Value = 24
$ws.Range("K7").Value = 24
$ws.Range("M7").Value = 89

$ws.Range("H16").Value = 635.25
$ws.Range("J16").Value = 834.4
$ws.Range("L16").Value = 834.4
$ws.Range("N16").Value = -1408.4

$ws.Range("H31").Value = 3732.375
$ws.Range("I31").Value = 1156.875
$ws.Range("J31").Value = 6307.875
$ws.Range("K31").Value = 1156.875
$ws.Range("L31").Value = 6307.875
$ws.Range("M31").Value = -861.875
$ws.Range("N31").Value = -6897.875

$ws.Range("H34").Value = 3732.375
$ws.Range("I34").Value = 1156.875
$ws.Range("J34").Value = 6307.875
$ws.Range("K34").Value = 1156.875
$ws.Range("L34").Value = 6307.875
$ws.Range("M34").Value = -954.875
$ws.Range("N34").Value = -6711.875

$ws.Range("H58").Value = 1518.96
$ws.Range("I58").Value = 741.4706
$ws.Range("K58").Value = 741.4706
$ws.Range("M58").Value = -538.4706

$ws.Range("H86").Value = 61101.715
$ws.Range("I86").Value = 89993.75
$ws.Range("J86").Value = 22579
$ws.Range("K86").Value = 89993.75
$ws.Range("L86").Value = 22579
$ws.Range("M86").Value = -88870.75
$ws.Range("N86").Value = -24825

$ws.Range("H89").Value = 61101.715
$ws.Range("I89").Value = 89993.75
$ws.Range("J89").Value = 22579
$ws.Range("K89").Value = 449968.75
$ws.Range("L89").Value = 112895
$ws.Range("M89").Value = -444352.75
$ws.Range("N89").Value = -124127

$ws.Range("H99").Value = 13206044
$ws.Range("I99").Value = 3489963.2
$ws.Range("K99").Value = 3489963.2
$ws.Range("M99").Value = -3488465.2

$ws.Range("H105").Value = 13894674
$ws.Range("I105").Value = 1651.9286
$ws.Range("J105").Value = 62520250
$ws.Range("K105").Value = 1651.9286
$ws.Range("L105").Value = 62520250
$ws.Range("M105").Value = 95.07140000000004
$ws.Range("N105").Value = -62523744

$ws.Range("H113").Value = 635.25
$ws.Range("J113").Value = 834.4
$ws.Range("L113").Value = 834.4
$ws.Range("N113").Value = -5174.4

$ws.Range("H120").Value = 689372.25
$ws.Range("J120").Value = 689372.25
$ws.Range("L120").Value = 689372.25
$ws.Range("N120").Value = -696630.25

$ws.Range("H122").Value = 333307.3
$ws.Range("I122").Value = 639546.75
$ws.Range("K122").Value = 1918640.25
$ws.Range("M122").Value = -1916190.25

$ws.Range("H126").Value = 13206044
$ws.Range("I126").Value = 3489963.2
$ws.Range("K126").Value = 10469889.6
$ws.Range("M126").Value = -10467419.6

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 1787.2858
$ws.Range("I132").Value = 1501.6923
$ws.Range("K132").Value = 4505.0769
$ws.Range("M132").Value = -1975.0769

$ws.Range("H134").Value = 2918.9343
$ws.Range("J134").Value = 6406.4375
$ws.Range("L134").Value = 19219.3125
$ws.Range("N134").Value = -24289.3125

$ws.Range("H136").Value = 1518.96
$ws.Range("I136").Value = 741.4706
$ws.Range("K136").Value = 2224.4118
$ws.Range("M136").Value = 325.5882000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6646.5557
$ws.Range("I113").Value = 34774.332
$ws.Range("J113").Value = 1021
$ws.Range("K113").Value = 104322.996
$ws.Range("L113").Value = 3063
$ws.Range("M113").Value = -102152.996
$ws.Range("N113").Value = -7403

$ws.Range("H140").Value = 1421.2941
$ws.Range("I140").Value = 1421.2941
$ws.Range("K140").Value = 4263.8823
$ws.Range("M140").Value = 916.1176999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3514.1155
$ws.Range("I122").Value = 1514.25
$ws.Range("J122").Value = 5228.2856
$ws.Range("K122").Value = 4542.75
$ws.Range("L122").Value = 15684.8568
$ws.Range("M122").Value = -2092.75
$ws.Range("N122").Value = -20584.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11269.223
$ws.Range("I40").Value = 16648.715
$ws.Range("K40").Value = 16648.715
$ws.Range("M40").Value = -16512.715

$ws.Range("H46").Value = 1218
$ws.Range("I46").Value = 1411.1666
$ws.Range("J46").Value = 986.2
$ws.Range("K46").Value = 1411.1666
$ws.Range("L46").Value = 986.2
$ws.Range("M46").Value = -1223.1666
$ws.Range("N46").Value = -1362.2

$ws.Range("H61").Value = 2721.257
$ws.Range("I61").Value = 2370.9565
$ws.Range("K61").Value = 2370.9565
$ws.Range("M61").Value = -2168.9565

$ws.Range("H100").Value = 66219.336
$ws.Range("I100").Value = 161842.58
$ws.Range("J100").Value = 5368.1816
$ws.Range("K100").Value = 161842.58
$ws.Range("L100").Value = 5368.1816
$ws.Range("M100").Value = -161301.58
$ws.Range("N100").Value = -6450.1816

$ws.Range("H113").Value = 2721.257
$ws.Range("I113").Value = 2370.9565
$ws.Range("K113").Value = 2370.9565
$ws.Range("M113").Value = -200.9564999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 950.5
$ws.Range("I107").Value = 950.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2851.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -931.5
$ws.Range("N107").ClearContents()

$ws.Range("H136").Value = 1880
$ws.Range("I136").Value = 933.75
$ws.Range("J136").Value = 3351.9443
$ws.Range("K136").Value = 2801.25
$ws.Range("L136").Value = 10055.8329
$ws.Range("M136").Value = -251.25
$ws.Range("N136").Value = -15155.8329
